$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2113.913
$ws.Range("I40").Value = 1904.4445
$ws.Range("J40").Value = 2248.5715
$ws.Range("K40").Value = 1904.4445
$ws.Range("L40").Value = 2248.5715
$ws.Range("M40").Value = -1729.4445
$ws.Range("N40").Value = -2598.5715
$ws.Range("H70").Value = 1743.3
$ws.Range("I70").Value = 1200
$ws.Range("J70").Value = 1976.1428
$ws.Range("K70").Value = 3600
$ws.Range("L70").Value = 5928.428400000001
$ws.Range("M70").Value = -3330
$ws.Range("N70").Value = -6468.428400000001
$ws.Range("H73").Value = 1743.3
$ws.Range("I73").Value = 1200
$ws.Range("J73").Value = 1976.1428
$ws.Range("K73").Value = 3600
$ws.Range("L73").Value = 5928.428400000001
$ws.Range("M73").Value = -2664
$ws.Range("N73").Value = -7800.428400000001
$ws.Range("H74").Value = 6893.6665
$ws.Range("J74").Value = 7536.4443
$ws.Range("L74").Value = 7536.4443
$ws.Range("N74").Value = -9408.444299999999
$ws.Range("H77").Value = 6893.6665
$ws.Range("J77").Value = 7536.4443
$ws.Range("L77").Value = 37682.2215
$ws.Range("N77").Value = -47042.2215
$ws.Range("H92").Value = 1045.3077
$ws.Range("I92").Value = 303
$ws.Range("J92").Value = 2233
$ws.Range("K92").Value = 303
$ws.Range("L92").Value = 2233
$ws.Range("M92").Value = 945
$ws.Range("N92").Value = -4729
$ws.Range("H99").Value = 1797.4286
$ws.Range("I99").Value = 257
$ws.Range("J99").Value = 2054.1667
$ws.Range("K99").Value = 771
$ws.Range("L99").Value = 6162.500100000001
$ws.Range("M99").Value = 727
$ws.Range("N99").Value = -9158.500100000001
$ws.Range("H113").Value = 3835.7778
$ws.Range("I113").Value = 3588.6
$ws.Range("K113").Value = 3588.6
$ws.Range("M113").Value = -334.5999999999999
$ws.Range("H132").Value = 82858.414
$ws.Range("I132").Value = 88069.96000000001
$ws.Range("K132").Value = 264209.88
$ws.Range("M132").Value = -261679.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27901.268
$ws.Range("I32").Value = 30078.658
$ws.Range("K32").Value = 30078.658
$ws.Range("M32").Value = -29791.658
$ws.Range("H45").Value = 1869.2858
$ws.Range("I45").Value = 2060.4167
$ws.Range("J45").Value = 722.5
$ws.Range("K45").Value = 2060.4167
$ws.Range("L45").Value = 722.5
$ws.Range("M45").Value = -1683.4167
$ws.Range("N45").Value = -1476.5
$ws.Range("H74").Value = 2011.7587
$ws.Range("I74").Value = 834.9545000000001
$ws.Range("K74").Value = 834.9545000000001
$ws.Range("M74").Value = 39.04549999999995
$ws.Range("H77").Value = 2011.7587
$ws.Range("I77").Value = 834.9545000000001
$ws.Range("K77").Value = 4174.7725
$ws.Range("M77").Value = 193.2275
$ws.Range("H102").Value = 26904.143
$ws.Range("I102").Value = 28919.46
$ws.Range("K102").Value = 28919.46
$ws.Range("M102").Value = -27297.46
$ws.Range("H122").Value = 1272.45
$ws.Range("I122").Value = 1263.9333
$ws.Range("K122").Value = 3791.7999
$ws.Range("M122").Value = -1341.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1917.8334
$ws.Range("J20").Value = 1966.3334
$ws.Range("L20").Value = 1966.3334
$ws.Range("N20").Value = -2460.3334
$ws.Range("H94").Value = 1144.0714
$ws.Range("I94").Value = 891.5454999999999
$ws.Range("K94").Value = 891.5454999999999
$ws.Range("M94").Value = -440.5454999999999
$ws.Range("H99").Value = 11503.23
$ws.Range("I99").Value = 4747.2856
$ws.Range("J99").Value = 19385.166
$ws.Range("K99").Value = 4747.2856
$ws.Range("L99").Value = 19385.166
$ws.Range("M99").Value = -3249.2856
$ws.Range("N99").Value = -22381.166
$ws.Range("H128").Value = 14473.333
$ws.Range("I128").Value = 14473.333
$ws.Range("K128").Value = 43419.999
$ws.Range("M128").Value = -40929.999
$ws.Range("H134").Value = 605895.7
$ws.Range("I134").Value = 568491.4399999999
$ws.Range("K134").Value = 1705474.32
$ws.Range("M134").Value = -1702939.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 101632.71
$ws.Range("I86").Value = 6062
$ws.Range("K86").Value = 6062
$ws.Range("M86").Value = -4939
$ws.Range("H89").Value = 101632.71
$ws.Range("I89").Value = 6062
$ws.Range("K89").Value = 30310
$ws.Range("M89").Value = -24694
$ws.Range("H107").Value = 927.1852
$ws.Range("I107").Value = 738
$ws.Range("K107").Value = 738
$ws.Range("M107").Value = 1182
$ws.Range("H120").Value = 23485.8
$ws.Range("J120").Value = 23485.8
$ws.Range("L120").Value = 23485.8
$ws.Range("N120").Value = -30743.8
$ws.Range("H122").Value = 3144.2083
$ws.Range("I122").Value = 2722.1177
$ws.Range("K122").Value = 8166.353099999999
$ws.Range("M122").Value = -5716.353099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 468
$ws.Range("I5").Value = 255.4
$ws.Range("J5").Value = 999.5
$ws.Range("K5").Value = 766.2
$ws.Range("L5").Value = 2998.5
$ws.Range("M5").Value = -654.2
$ws.Range("N5").Value = -3222.5
$ws.Range("H131").Value = 9526.143
$ws.Range("J131").Value = 11360.087
$ws.Range("L131").Value = 34080.261
$ws.Range("N131").Value = -44160.261
$ws.Range("H135").Value = 468
$ws.Range("I135").Value = 255.4
$ws.Range("J135").Value = 999.5
$ws.Range("K135").Value = 2298.6
$ws.Range("L135").Value = 8995.5
$ws.Range("M135").Value = 236.4000000000001
$ws.Range("N135").Value = -14065.5
$ws.Range("H137").Value = 2251.3125
$ws.Range("I137").Value = 1152.091
$ws.Range("J137").Value = 4669.6
$ws.Range("K137").Value = 3456.273
$ws.Range("L137").Value = 14008.8
$ws.Range("M137").Value = 1643.727
$ws.Range("N137").Value = -24208.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5311.3335
$ws.Range("I70").Value = 5426.3335
$ws.Range("K70").Value = 5426.3335
$ws.Range("M70").Value = -5156.3335
$ws.Range("H73").Value = 5311.3335
$ws.Range("I73").Value = 5426.3335
$ws.Range("K73").Value = 5426.3335
$ws.Range("M73").Value = -4490.3335
$ws.Range("H97").Value = 1207.7307
$ws.Range("I97").Value = 1394.35
$ws.Range("K97").Value = 1394.35
$ws.Range("M97").Value = -898.3499999999999
$ws.Range("H102").Value = 2544.0557
$ws.Range("I102").Value = 1916.8966
$ws.Range("J102").Value = 5142.2856
$ws.Range("K102").Value = 1916.8966
$ws.Range("L102").Value = 5142.2856
$ws.Range("M102").Value = -294.8966
$ws.Range("N102").Value = -8386.285599999999
$ws.Range("H123").Value = 74999.336
$ws.Range("J123").Value = 74999.336
$ws.Range("L123").Value = 74999.336
$ws.Range("N123").Value = -79899.336
$ws.Range("H132").Value = 451381.78
$ws.Range("I132").Value = 552128.75
$ws.Range("J132").Value = 8095.2
$ws.Range("K132").Value = 1656386.25
$ws.Range("L132").Value = 24285.6
$ws.Range("M132").Value = -1653856.25
$ws.Range("N132").Value = -29345.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10113.728
$ws.Range("J16").Value = 1193.7142
$ws.Range("L16").Value = 1193.7142
$ws.Range("N16").Value = -1533.7142
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()  # remove N24 (was -3061)
$ws.Range("H122").Value = 5161.1514
$ws.Range("I122").Value = 4573.08
$ws.Range("K122").Value = 13719.24
$ws.Range("M122").Value = -11269.24
$ws.Range("H132").Value = 1028028.4
$ws.Range("I132").Value = 1292661.8
$ws.Range("K132").Value = 3877985.4
$ws.Range("M132").Value = -3875455.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2726.7778
$ws.Range("J96").Value = 4566
$ws.Range("L96").Value = 4566
$ws.Range("N96").Value = -7312
$ws.Range("H100").Value = 1884.4546
$ws.Range("I100").Value = 1420.2858
$ws.Range("K100").Value = 2840.5716
$ws.Range("M100").Value = -2299.5716
$ws.Range("H126").Value = 4174.04
$ws.Range("I126").Value = 4006.8635
$ws.Range("K126").Value = 12020.5905
$ws.Range("M126").Value = -9550.5905
